$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct writes. D-column numeric-looking strings are written
# with a leading apostrophe to force text, then the style is reset to
# "Normal" so no stray NumberFormat/quote-prefix style sticks on the cell
# (matches the source file, where these cells carry no explicit style).

$ws.Range("D2").Value = "'54.045.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "'2.241.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'494.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").Value = "'127.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("D7").Value = "'0.993"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").Value = "'2.281.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("D10").Value = "'0.0948"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("D12").Value = "'0.325"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.50%  "
$ws.Range("D13").Value = "'4.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").Value = "'2.646.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "'21.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.79%  "
$ws.Range("D16").Value = "'53.951.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "'2.275.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").Value = "'10.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.30%  "
$ws.Range("E20").Value = "  +3.63%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'301.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'6.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.89%  "
$ws.Range("D23").Value = "'0.994"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  +2.77%  "
$ws.Range("D28").Value = "'2.371.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("D30").Value = "'7.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").Value = "'168.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("E34").Value = "  +2.72%  "
$ws.Range("D35").Value = "'0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'0.991"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").Value = "'17.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("D40").Value = "'0.864"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.55%  "
$ws.Range("D42").Value = "'35.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  +3.55%  "
$ws.Range("E44").Value = "  +2.82%  "
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'4.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.85%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'127.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.36%  "
$ws.Range("D48").Value = "'0.0889"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("D49").Value = "'0.542"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").Value = "'237.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.92%  "
$ws.Range("E51").Value = "  +3.25%  "
